$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Ben Stokes"
$ws.Name = "Ben Stokes"

# Insert a new column at A, shifting teamName..result (old A:L) right to B:M
$ws.Range("A1").EntireColumn.Insert()

# Populate the new "matchNo" column
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "4th"
